$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 129: mentions ---
$cA129 = $ws.Cells.Item(129, 1)
$ws.Hyperlinks.Add($cA129, "http://purl.obolibrary.org/obo/IAO_0000142") | Out-Null
$ws.Cells.Item(129, 2).Value = "mentions"
$ws.Cells.Item(129, 3).Value = "y"

# --- Row 130: plan ---
$ws.Cells.Item(130, 2).Value = "plan"
$cA130 = $ws.Cells.Item(130, 1)
$ws.Hyperlinks.Add($cA130, "http://purl.obolibrary.org/obo/OBI_0000260") | Out-Null
$ws.Cells.Item(130, 3).Value = "y"

# --- Row 131: categorical value specification ---
$cB131 = $ws.Cells.Item(131, 2)
$cB131.Value = "categorical value specification"
$cA131 = $ws.Cells.Item(131, 1)
$ws.Hyperlinks.Add($cA131, "http://purl.obolibrary.org/obo/OBI_0001930") | Out-Null
$ws.Hyperlinks.Add($cB131, "http://purl.obolibrary.org/obo/OBI_0001930", "", "", "http://purl.obolibrary.org/obo/OBI_0001930 categorical value specification") | Out-Null
$cB131.Value = "categorical value specification"
$ws.Cells.Item(131, 3).Value = "y"
